# 1. Slide 5's table (the financial documents comparison table) gets a new
#    built-in table style applied (PowerPoint "Table Design" gallery pick).
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$sh = $s.Shapes.Item(2)
$tbl = $sh.Table
$tbl.ApplyStyle("{6D5853ED-1638-4FB8-8072-E98553343E78}")

# 2. The deck's theme (linked from the slide master / theme1.xml) is switched
#    from the custom "Integral" (Red Violet) palette to the standard Office
#    theme palette. dk1/lt1 (black/white) are unchanged; the other ten theme
#    colour slots move to the stock Office values.
function HexToOle($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$officeColors = @{
    3  = "44546A"   # dk2
    4  = "E7E6E6"   # lt2
    5  = "5B9BD5"   # accent1
    6  = "ED7D31"   # accent2
    7  = "A5A5A5"   # accent3
    8  = "FFC000"   # accent4
    9  = "4472C4"   # accent5
    10 = "70AD47"   # accent6
    11 = "0563C1"   # hlink
    12 = "954F72"   # folHlink
}

$cs = $p.Slides.Item(1).ThemeColorScheme
foreach ($idx in $officeColors.Keys) {
    $cs.Item($idx).RGB = HexToOle $officeColors[$idx]
}
